$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table was regenerated so that the "Enterprises density
# (per 1000 people)" / "69.3" pair now comes right after "MSMEs", ahead of
# "Employment (% of total)" / "32" and "Enterprises (absolute #)" / "17000".
# Because the cells in column A/D keep referencing the same string-table
# slots, the net visible effect is a 3-row rotation of the label/value
# pairs that sit in rows 12-14:
#   row 12 (was Employment / 32)            -> Enterprises density / 69.3
#   row 13 (was Enterprises (absolute#) / 17000) -> Employment / 32
#   row 14 (was Enterprises density / 69.3)  -> Enterprises (absolute #) / 17000
#
# The numeric-looking values are stored as text in the workbook, so force
# text formatting before assigning them, then restore the default "Normal"
# style so no stray number formatting is left behind.

$ws.Range("A12").Value = "Enterprises density (per 1000 people)"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "69.3"
$ws.Range("D12").Style = "Normal"

$ws.Range("A13").Value = "Employment (% of total)"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32"
$ws.Range("D13").Style = "Normal"

$ws.Range("A14").Value = "Enterprises (absolute #)"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17000"
$ws.Range("D14").Style = "Normal"
